$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C from 45184 -> 45186
# for every data row that currently holds that value.
for ($r = 2; $r -le 176; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $val = $cCell.Value2
    if ($val -eq 45184) {
        $cCell.Value = 45186
    }
}

# Add a friendly-name second argument to every HYPERLINK() formula in
# columns S..Y (Artfyndslänk .. Tillsynsbegäransmaillänk). The friendly
# name is the row's "Beteckning" value in column A, which matches the
# file-name stem already used inside the URL.
$cols = 19..25  # S=19, T=20, U=21, V=22, W=23, X=24, Y=25

for ($r = 2; $r -le 176; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($name)) { continue }

    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -match '^=HYPERLINK\("([^"]*)"\)$') {
                $url = $matches[1]
                $cell.Formula = '=HYPERLINK("' + $url + '", "' + $name + '")'
            }
        }
    }
}
